$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "2025-08-27T08:00"
$ws.Cells.Item(2, 2).Value = 21.6
$ws.Cells.Item(2, 3).Value = 77
$ws.Cells.Item(3, 1).Value = "2025-08-27T09:00"
$ws.Cells.Item(3, 2).Value = 23.1
$ws.Cells.Item(3, 3).Value = 67
$ws.Cells.Item(4, 1).Value = "2025-08-27T10:00"
$ws.Cells.Item(4, 2).Value = 24.5
$ws.Cells.Item(4, 3).Value = 60
$ws.Cells.Item(5, 1).Value = "2025-08-27T11:00"
$ws.Cells.Item(5, 2).Value = 26
$ws.Cells.Item(5, 3).Value = 57
$ws.Cells.Item(6, 1).Value = "2025-08-27T12:00"
$ws.Cells.Item(6, 2).Value = 26.7
$ws.Cells.Item(6, 3).Value = 48
$ws.Cells.Item(7, 1).Value = "2025-08-27T13:00"
$ws.Cells.Item(7, 2).Value = 26.9
$ws.Cells.Item(7, 3).Value = 48
$ws.Cells.Item(8, 1).Value = "2025-08-27T14:00"
$ws.Cells.Item(8, 2).Value = 26.2
$ws.Cells.Item(8, 3).Value = 52
$ws.Cells.Item(9, 1).Value = "2025-08-27T15:00"
$ws.Cells.Item(9, 2).Value = 24.2
$ws.Cells.Item(9, 3).Value = 68
$ws.Cells.Item(10, 1).Value = "2025-08-27T16:00"
$ws.Cells.Item(10, 2).Value = 23
$ws.Cells.Item(10, 3).Value = 72
$ws.Cells.Item(11, 1).Value = "2025-08-27T17:00"
$ws.Cells.Item(11, 2).Value = 22.5
$ws.Cells.Item(11, 3).Value = 62
$ws.Cells.Item(12, 1).Value = "2025-08-27T18:00"
$ws.Cells.Item(12, 2).Value = 20.8
$ws.Cells.Item(12, 3).Value = 82
$ws.Cells.Item(13, 1).Value = "2025-08-27T19:00"
$ws.Cells.Item(13, 2).Value = 20.3
$ws.Cells.Item(13, 3).Value = 80
$ws.Cells.Item(14, 1).Value = "2025-08-27T20:00"
$ws.Cells.Item(14, 2).Value = 19.8
$ws.Cells.Item(14, 3).Value = 81
$ws.Cells.Item(15, 1).Value = "2025-08-27T21:00"
$ws.Cells.Item(15, 2).Value = 19.3
$ws.Cells.Item(15, 3).Value = 87
$ws.Cells.Item(16, 1).Value = "2025-08-27T22:00"
$ws.Cells.Item(16, 2).Value = 18.9
$ws.Cells.Item(16, 3).Value = 92
$ws.Cells.Item(17, 1).Value = "2025-08-27T23:00"
$ws.Cells.Item(17, 2).Value = 18.6
$ws.Cells.Item(17, 3).Value = 92
$ws.Cells.Item(18, 1).Value = "2025-08-28T00:00"
$ws.Cells.Item(18, 2).Value = 18
$ws.Cells.Item(18, 3).Value = 97
$ws.Cells.Item(19, 1).Value = "2025-08-28T01:00"
$ws.Cells.Item(19, 2).Value = 18.3
$ws.Cells.Item(19, 3).Value = 87
$ws.Cells.Item(20, 1).Value = "2025-08-28T02:00"
$ws.Cells.Item(20, 2).Value = 17.9
$ws.Cells.Item(20, 3).Value = 87
$ws.Cells.Item(21, 1).Value = "2025-08-28T03:00"
$ws.Cells.Item(21, 2).Value = 17.8
$ws.Cells.Item(21, 3).Value = 90
$ws.Cells.Item(22, 1).Value = "2025-08-28T04:00"
$ws.Cells.Item(22, 2).Value = 17.9
$ws.Cells.Item(22, 3).Value = 90
$ws.Cells.Item(23, 1).Value = "2025-08-28T05:00"
$ws.Cells.Item(23, 2).Value = 17.4
$ws.Cells.Item(23, 3).Value = 94
$ws.Cells.Item(24, 1).Value = "2025-08-28T06:00"
$ws.Cells.Item(24, 2).Value = 17.2
$ws.Cells.Item(24, 3).Value = 94
$ws.Cells.Item(25, 1).Value = "2025-08-28T07:00"
$ws.Cells.Item(25, 2).Value = 17.3
$ws.Cells.Item(25, 3).Value = 92
$ws.Cells.Item(26, 1).Value = "2025-08-28T08:00"
$ws.Cells.Item(26, 2).Value = 17.5
$ws.Cells.Item(26, 3).Value = 88
$ws.Cells.Item(27, 1).Value = "2025-08-28T09:00"
$ws.Cells.Item(27, 2).Value = 17.4
$ws.Cells.Item(27, 3).Value = 90
$ws.Cells.Item(28, 1).Value = "2025-08-28T10:00"
$ws.Cells.Item(28, 2).Value = 17.4
$ws.Cells.Item(28, 3).Value = 85
$ws.Cells.Item(29, 1).Value = "2025-08-28T11:00"
$ws.Cells.Item(29, 2).Value = 18.1
$ws.Cells.Item(29, 3).Value = 80
$ws.Cells.Item(30, 1).Value = "2025-08-28T12:00"
$ws.Cells.Item(30, 2).Value = 18.5
$ws.Cells.Item(30, 3).Value = 72
$ws.Cells.Item(31, 1).Value = "2025-08-28T13:00"
$ws.Cells.Item(31, 2).Value = 18.6
$ws.Cells.Item(31, 3).Value = 68
$ws.Cells.Item(32, 1).Value = "2025-08-28T14:00"
$ws.Cells.Item(32, 2).Value = 18.7
$ws.Cells.Item(32, 3).Value = 69
$ws.Cells.Item(33, 1).Value = "2025-08-28T15:00"
$ws.Cells.Item(33, 2).Value = 18.6
$ws.Cells.Item(33, 3).Value = 78
$ws.Cells.Item(34, 1).Value = "2025-08-28T16:00"
$ws.Cells.Item(34, 2).Value = 18.6
$ws.Cells.Item(34, 3).Value = 72
$ws.Cells.Item(35, 1).Value = "2025-08-28T17:00"
$ws.Cells.Item(35, 2).Value = 18.4
$ws.Cells.Item(35, 3).Value = 76
$ws.Cells.Item(36, 1).Value = "2025-08-28T18:00"
$ws.Cells.Item(36, 2).Value = 18.5
$ws.Cells.Item(36, 3).Value = 75
$ws.Cells.Item(37, 1).Value = "2025-08-28T19:00"
$ws.Cells.Item(37, 2).Value = 17.9
$ws.Cells.Item(37, 3).Value = 79
$ws.Cells.Item(38, 1).Value = "2025-08-28T20:00"
$ws.Cells.Item(38, 2).Value = 17.2
$ws.Cells.Item(38, 3).Value = 85
$ws.Cells.Item(39, 1).Value = "2025-08-28T21:00"
$ws.Cells.Item(39, 2).Value = 17.2
$ws.Cells.Item(39, 3).Value = 87
$ws.Cells.Item(40, 1).Value = "2025-08-28T22:00"
$ws.Cells.Item(40, 2).Value = 16.9
$ws.Cells.Item(40, 3).Value = 92
$ws.Cells.Item(41, 1).Value = "2025-08-28T23:00"
$ws.Cells.Item(41, 2).Value = 16.7
$ws.Cells.Item(41, 3).Value = 94
$ws.Cells.Item(42, 1).Value = "2025-08-29T00:00"
$ws.Cells.Item(42, 2).Value = 15.5
$ws.Cells.Item(42, 3).Value = 89
$ws.Cells.Item(43, 1).Value = "2025-08-29T01:00"
$ws.Cells.Item(43, 2).Value = 15
$ws.Cells.Item(43, 3).Value = 90
$ws.Cells.Item(44, 1).Value = "2025-08-29T02:00"
$ws.Cells.Item(44, 2).Value = 14.8
$ws.Cells.Item(44, 3).Value = 89
$ws.Cells.Item(45, 1).Value = "2025-08-29T03:00"
$ws.Cells.Item(45, 2).Value = 14.6
$ws.Cells.Item(45, 3).Value = 88
$ws.Cells.Item(46, 1).Value = "2025-08-29T04:00"
$ws.Cells.Item(46, 2).Value = 14.5
$ws.Cells.Item(46, 3).Value = 84
$ws.Cells.Item(47, 1).Value = "2025-08-29T05:00"
$ws.Cells.Item(47, 2).Value = 14.3
$ws.Cells.Item(47, 3).Value = 85
$ws.Cells.Item(48, 1).Value = "2025-08-29T06:00"
$ws.Cells.Item(48, 2).Value = 15
$ws.Cells.Item(48, 3).Value = 85
$ws.Cells.Item(49, 1).Value = "2025-08-29T07:00"
$ws.Cells.Item(49, 2).Value = 15.9
$ws.Cells.Item(49, 3).Value = 78

$wb.Save()
